$wb = $excel.ActiveWorkbook

$wsPre = $wb.Worksheets.Item("RegrasEmailsPreTratamento")
$wsPre.Activate()
$wsPre.Range("A5").Select()

$wsDiscard = $wb.Worksheets.Item("RegrasEmailDiscard")
$wsDiscard.Range("E2").Value = "Tratamento Manual"
$wsDiscard.Range("E3").Value = "Tratamento Manual"
$wsDiscard.Activate()
$wsDiscard.Range("A7").Select()
